# Rename the "ac_w" variable to "x_ac" and update its description
# from "Aerodynamic center from wing" to "Aerodynamic center from nose".
# This corresponds to row 25 on the "Data" worksheet:
#   B25: ac_w -> x_ac
#   E25: Aerodynamic center from wing -> Aerodynamic center from nose

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("E25").Value = "Aerodynamic center from nose"
$ws.Range("B25").Value = "x_ac"

# Reflect the author's cursor position at the time of the edit.
$ws.Range("C23").Select()
